$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.992.43"
$ws.Range("E2").Value = "  +0.65%  "

# Row 3
$ws.Range("D3").Value = "2.039.78"
$ws.Range("E3").Value = "  -3.80%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.63"
$ws.Range("E5").Value = "  -3.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.647"
$ws.Range("E6").Value = "  -3.42%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.18"
$ws.Range("E8").Value = "  +14.73%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.05"
$ws.Range("E9").Value = "  -0.21%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.372"
$ws.Range("E10").Value = "  -0.84%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0749"
$ws.Range("E11").Value = "  +1.01%  "

# Row 12
$ws.Range("E12").Value = "  +4.95%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.83"
$ws.Range("E13").Value = "  +2.37%  "

# Row 14
$ws.Range("D14").Value = "2.342.17"
$ws.Range("E14").Value = "  -3.56%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.804"
$ws.Range("E15").Value = "  -5.05%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.044.80"
$ws.Range("E16").Value = "  -3.59%  "

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.12"
$ws.Range("E17").Value = "  -0.36%  "

# Row 18
$ws.Range("D18").Value = "36.916.64"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.26"
$ws.Range("E19").Value = "  -3.51%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0883"
$ws.Range("E20").Value = "  +4.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.05"
$ws.Range("E21").Value = "  +4.96%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "234.58"
$ws.Range("E22").Value = "  -3.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.17"
$ws.Range("E23").Value = "  -1.35%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("E25").Value = "  -4.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.01"
$ws.Range("E26").Value = "  -2.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.91"
$ws.Range("E27").Value = "  -3.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.72"
$ws.Range("E28").Value = "  -8.53%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.96"
$ws.Range("E29").Value = "  -3.34%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.121"
$ws.Range("E30").Value = "  -1.95%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.47"
$ws.Range("E31").Value = "  -1.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.04"
$ws.Range("E32").Value = "  +9.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0610"
$ws.Range("E33").Value = "  +1.27%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.26"
$ws.Range("E34").Value = "  +1.64%  "

# Row 35
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.28%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0869"
$ws.Range("E36").Value = "  -9.11%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.22"
$ws.Range("E37").Value = "  -5.80%  "

# Row 38
$ws.Range("E38").Value = "  -6.30%  "

# Row 39
$ws.Range("E39").Value = "  -3.32%  "

# Row 40
$ws.Range("E40").Value = "  +20.62%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.97"
$ws.Range("E41").Value = "  +10.34%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0220"
$ws.Range("E42").Value = "  -2.65%  "

# Row 43
$ws.Range("B43").Value = "Gas"
$ws.Range("C43").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.08"
$ws.Range("E43").Value = "  -44.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.12"
$ws.Range("E44").Value = "  -6.71%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.54"
$ws.Range("E45").Value = "  -4.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.77"
$ws.Range("E46").Value = "  -1.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.96"
$ws.Range("E47").Value = "  +37.61%  "

# Row 48
$ws.Range("D48").Value = "1.281.26"
$ws.Range("E48").Value = "  -5.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  +2.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  +0.78%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.69"
$ws.Range("E51").Value = "  -6.65%  "
